$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,4,5,6,7,10,11,13)

$data = @(
    @(0.1422219981393198,0.0236562874145676,0.06297198839646345,5.93924082002718,0.002640430703189177,0.2049883900821001,1.58633087707787,0.4352249282588616),
    @(0.1332389306771375,0.02111473659372365,0.06350821472012491,5.756369236366282,0.002646582247319672,0.2018037762405527,1.551979390455529,0.4313032686060936),
    @(0.1278006476182725,0.01956116484976889,0.06386878086799364,5.645234355723602,0.002650553865668184,0.199895959969659,1.532595508604857,0.4293178096188228),
    @(0.1256039844385981,0.01892967298528703,0.06402360523017681,5.600229505118534,0.002652221433860459,0.199130373132661,1.525124112823136,0.4286146895431528),
    @(0.1252404077234104,0.01882490668131709,0.06404979071728079,5.592773450594649,0.002652501302927458,0.1990039618791073,1.523909272695647,0.428504329942939),
    @(0.1277709437305248,0.01955264203804319,0.06387083691461015,5.644626263790656,0.002650576155940848,0.199885587081198,1.532493017430824,0.429307898368414),
    @(0.1391085951186426,0.0227784290155455,0.06315038735699829,5.875944890208189,0.002642511487363905,0.2038804088939514,1.574130851720724,0.4337848973835534),
    @(0.1619555848001539,0.02916622327664697,0.06198547473681515,6.338946363880041,0.002628232026333055,0.2120962445639947,1.669426863508875,0.4459298206112976),
    @(0.1791177055421542,0.0339070825373966,0.06127991148028755,6.685245443518369,0.002618665235652102,0.2183725938110896,1.747899756570746,0.4569258596651409),
    @(0.1870075621961433,0.03607636725023156,0.06099140135533077,6.844208980118992,0.00261451126866552,0.2212816639493624,1.785468009893009,0.4623833745018615),
    @(0.1900071510025754,0.03689980724732322,0.06088680378903177,6.90461683841022,0.002612966551834872,0.2223911292819878,1.799965581916183,0.4645158539777441),
    @(0.1893606082664547,0.03672237412041568,0.06090912390236625,6.891597405935784,0.00261329797842504,0.2221518345797762,1.796831167256528,0.4640536524368599),
    @(0.1872541023845997,0.03614407156641164,0.06098270282556761,6.849174490364135,0.002614383617752172,0.2213727821001044,1.786655283369271,0.4625574933133905),
    @(0.1859653520231745,0.03579010707031216,0.06102837790490412,6.823216998190674,0.002615052283450916,0.2208966173039073,1.780457660199772,0.4616496373683674),
    @(0.1786037473397926,0.03376558281765085,0.06129941833083308,6.674886107005733,0.002618940677038052,0.2181835736562689,1.745482418807427,0.4565783914815498),
    @(0.1741088107454232,0.03252695510667536,0.06147399656702746,6.584260590507938,0.002621376672337151,0.2165330981926772,1.724506953883491,0.4535842253022508),
    @(0.1715312295997933,0.03181570962253488,0.06157746478234039,6.532269448865605,0.002622796441101434,0.2155888576220377,1.712618343534245,0.4519048988341225),
    @(0.1706598431823494,0.03157509198686625,0.0616130224829643,6.514689015308704,0.00262328035858641,0.2152700213562184,1.708623209842216,0.4513436546652585),
    @(0.1745864986220766,0.03265868581303266,0.06145509627759616,6.593893892078427,0.002621115427636502,0.2167082686915975,1.726721604127277,0.4538985232146757),
    @(0.1878725122224978,0.03631387804719566,0.06096096467144996,6.861629330985068,0.002614063972426093,0.2216013944037485,1.78963680569251,0.4629951618137795),
    @(0.1966248754483928,0.03871436808326223,0.06066514843950621,7.03784746605254,0.00260962031443931,0.2248452220866994,1.832337928446577,0.4693242639454169),
    @(0.1919472514772025,0.03743206547738964,0.06082055283233956,6.943681242629395,0.002611976949792815,0.2231096949909386,1.809401966424502,0.4659110530270496),
    @(0.1743705153182873,0.03259912770033679,0.06146363143905909,6.589538332172907,0.002621233476446099,0.2166290596654221,1.725719829968426,0.4537562981491376),
    @(0.1557089270916947,0.02743049289203725,0.06227416235906169,6.212645293008563,0.00263193183729944,0.2098320000983591,1.642171689717429,0.4386070558177845)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowVals[$j]
    }
}

Write-Output "done"
